$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 461, shifting the existing rows (461..565) down to (462..566)
$ws.Rows(461).Insert()

# Populate the newly inserted row 461 with its data
$ws.Cells.Item(461, 1).Value = 6
$ws.Cells.Item(461, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(461, 3).Value = "Metropolitana"
$ws.Cells.Item(461, 4).Value = 44889
$ws.Cells.Item(461, 5).Value = 13
$ws.Cells.Item(461, 6).Value = 100112039
$ws.Cells.Item(461, 7).Value = "Ciboulette"
$ws.Cells.Item(461, 8).Value = "Sin especificar"
$ws.Cells.Item(461, 9).Value = "Primera"
$ws.Cells.Item(461, 10).Value = 370
$ws.Cells.Item(461, 11).Value = 900
$ws.Cells.Item(461, 12).Value = 900
$ws.Cells.Item(461, 13).Value = 900
$ws.Cells.Item(461, 14).Value = '$/docena de atados'
$ws.Cells.Item(461, 15).Value = "Región Metropolitana"
$ws.Cells.Item(461, 16).Value = 300
$ws.Cells.Item(461, 17).Value = 3
$ws.Cells.Item(461, 18).Value = "Hortaliza"
